$d = $word.ActiveDocument

# Target the first paragraph: "This is a Microsoft word document."
$p1 = $d.Paragraphs(1).Range

# Add two trailing spaces to the existing run's text (keeps existing formatting).
$p1.InsertAfter("  ")

# Remember the end of the paragraph's text (before the paragraph mark) so we
# can select exactly the new run we are about to append.
$afterSpacesEnd = $p1.End

# Append the new annotation text; it will initially share formatting with
# the run before it (since no rPr differs yet).
$p1.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch alternate)")

# Select just the newly-inserted annotation text (End - 1 to skip the
# trailing paragraph mark both sides account for) and color it dark red
# (RGB 192,0,0 -> hex C00000), which forces it into its own run with its
# own <w:rPr><w:color w:val="C00000"/></w:rPr>.
$newRun = $d.Range($afterSpacesEnd - 1, $p1.End - 1)
$newRun.Font.Color = 192
